# Update loading_percent values for Case_3_87 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 6.917169992021162
    "C2" = 6.806290309034241
    "D2" = 4.560284049730241
    "F2" = 21.68750424701895
    "G2" = 24.26825617323693
    "H2" = 13.19785320135738
    "K2" = 8.381888689455606
    "M2" = 19.97288045019196
    "O2" = 19.48192466483431
    "B3" = 6.639050935440197
    "C3" = 6.795415270960925
    "D3" = 4.478172084062789
    "F3" = 21.71259478648555
    "G3" = 24.30719758079749
    "H3" = 13.24037575981718
    "K3" = 8.215348006471942
    "M3" = 19.37585180230192
    "O3" = 19.54673329719562
    "B4" = 6.462751919658459
    "C4" = 6.789050318798203
    "D4" = 4.426248764577132
    "F4" = 21.73385448644175
    "G4" = 24.33940674518097
    "H4" = 13.26848876668649
    "K4" = 8.110137913184863
    "M4" = 19.00793810844169
    "O4" = 19.59066358829743
    "B5" = 6.38962232704481
    "C5" = 6.786537135112956
    "D5" = 4.404725287460335
    "F5" = 21.74398689699724
    "G5" = 24.35460957943739
    "H5" = 13.28044893827955
    "K5" = 8.066555654247253
    "M5" = 18.85792650730504
    "O5" = 19.60960350247923
    "B6" = 6.377404839228422
    "C6" = 6.786124749310584
    "D6" = 4.401129763166075
    "F6" = 21.74575798950898
    "G6" = 24.35725916334712
    "H6" = 13.2824653489227
    "K6" = 8.059277113918803
    "M6" = 18.83301973544559
    "O6" = 19.61281107833718
    "B7" = 6.461770728655673
    "C7" = 6.789016096199149
    "D7" = 4.425959946914896
    "F7" = 21.73398519281687
    "G7" = 24.33960337839829
    "H7" = 13.26864802555586
    "K7" = 8.109552968765033
    "M7" = 19.00591499693682
    "O7" = 19.59091481926686
    "B8" = 6.822481993686577
    "C8" = 6.802476991420743
    "D8" = 4.532292919248945
    "F8" = 21.69493879048071
    "G8" = 24.27995535804337
    "H8" = 13.21209903824632
    "K8" = 8.32509850611142
    "M8" = 19.7674495233114
    "O8" = 19.50341048284109
    "B9" = 7.481940355591249
    "C9" = 6.831262640533452
    "D9" = 4.728262003596139
    "F9" = 21.66492819940611
    "G9" = 24.22920419161006
    "H9" = 13.11710840280199
    "K9" = 8.722933300723781
    "M9" = 21.23982844995322
    "O9" = 19.36475155589453
    "B10" = 7.983069072389217
    "C10" = 6.853747483971909
    "D10" = 4.863841776517855
    "F10" = 21.67137028969308
    "G10" = 24.23270788850866
    "H10" = 13.05701423204285
    "K10" = 8.998314751500097
    "M10" = 22.29597179208946
    "O10" = 19.28310149965978
    "B11" = 8.206771477333318
    "C11" = 6.864240596720784
    "D11" = 4.923553309155041
    "F11" = 21.68049439738698
    "G11" = 24.24322072677658
    "H11" = 13.03178154251733
    "K11" = 9.119587781364668
    "M11" = 22.76850178008058
    "O11" = 19.25037722834395
    "B12" = 8.289686105981081
    "C12" = 6.86824989176986
    "D12" = 4.945870986780098
    "F12" = 21.68483923439694
    "G12" = 24.24848686336804
    "H12" = 13.02252928631489
    "K12" = 9.164910227473046
    "M12" = 22.94612397224047
    "O12" = 19.2386231024474
    "B13" = 8.271909122103658
    "C13" = 6.867384863338668
    "D13" = 4.941077713752655
    "F13" = 21.68386394788561
    "G13" = 24.24729552357722
    "H13" = 13.02450845421764
    "K13" = 9.155176350600041
    "M13" = 22.90793095173935
    "O13" = 19.24112615766343
    "B14" = 8.213629015197697
    "C14" = 6.8645697398179
    "D14" = 4.925395347888837
    "F14" = 21.68083402620757
    "G14" = 24.24362820903766
    "H14" = 13.03101428380601
    "K14" = 9.123328673252503
    "M14" = 22.78314216202013
    "O14" = 19.24939741235879
    "B15" = 8.177696298969803
    "C15" = 6.862849987346284
    "D15" = 4.91575086388815
    "F15" = 21.67909394074268
    "G15" = 24.24154928957893
    "H15" = 13.03503873376942
    "K15" = 9.103742044398977
    "M15" = 22.7065292745365
    "O15" = 19.25454693160598
    "B16" = 7.968199126984508
    "C16" = 6.853066865608362
    "D16" = 4.859898975164074
    "F16" = 21.67089858105373
    "G16" = 24.23220058047945
    "H16" = 13.0587056114277
    "K16" = 8.990306537430817
    "M16" = 22.26491609385507
    "O16" = 19.28532923457212
    "B17" = 7.836498477056689
    "C17" = 6.847131405598385
    "D17" = 4.825124239410237
    "F17" = 21.66745694509444
    "G17" = 24.22875215431339
    "H17" = 13.07376361430489
    "K17" = 8.919674147862802
    "M17" = 21.99184241457897
    "O17" = 19.30534669021458
    "B18" = 7.759585142770427
    "C18" = 6.843742532656467
    "D18" = 4.804938461808146
    "F18" = 21.66606048152606
    "G18" = 24.22760814716939
    "H18" = 13.08262264805069
    "K18" = 8.878673252312986
    "M18" = 21.8340421089606
    "O18" = 19.31727597519237
    "B19" = 7.733344587955944
    "C19" = 6.842599489332522
    "D19" = 4.798072612800079
    "F19" = 21.6656878199746
    "G19" = 24.22736487560834
    "H19" = 13.08565617714836
    "K19" = 8.864727465581728
    "M19" = 21.78049296830377
    "O19" = 19.32138635988986
    "B20" = 7.8506387138571
    "C20" = 6.847760669644748
    "D20" = 4.828845224870416
    "F20" = 21.66776297095271
    "G20" = 24.22903234298813
    "H20" = 13.07214016341039
    "K20" = 8.92723209251729
    "M20" = 22.02098903391805
    "O20" = 19.30317274893974
    "B21" = 8.230796186526627
    "C21" = 6.865395656947983
    "D21" = 4.930009698485078
    "F21" = 21.68169985174364
    "G21" = 24.24467049620224
    "H21" = 13.02909514402028
    "K21" = 9.132699616629539
    "M21" = 22.81983265493885
    "O21" = 19.24695061448469
    "B22" = 8.468773004461493
    "C22" = 6.877128634031708
    "D22" = 4.994409415517992
    "F22" = 21.69599346217299
    "G22" = 24.26238192805654
    "H22" = 13.00272787170531
    "K22" = 9.263470639137394
    "M22" = 23.33419052798705
    "O22" = 19.2139251113945
    "B23" = 8.342724297698052
    "C23" = 6.87084829291358
    "D23" = 4.960198732605446
    "F23" = 21.68789077172912
    "G23" = 24.25224310530993
    "H23" = 13.01663900456811
    "K23" = 9.194005275118467
    "M23" = 23.06042945325657
    "O23" = 19.23121038387439
    "B24" = 7.844249639303605
    "C24" = 6.847476105992971
    "D24" = 4.827163568623036
    "F24" = 21.66762280324765
    "G24" = 24.22890305804275
    "H24" = 13.07287349642959
    "K24" = 8.923816369086891
    "M24" = 22.00781434614303
    "O24" = 19.30415427670754
    "B25" = 7.309178567555392
    "C25" = 6.82323196024339
    "D25" = 4.676670812115024
    "F25" = 21.66804824165675
    "G25" = 24.23579563805901
    "H25" = 13.14110337875124
    "K25" = 8.618157827501676
    "M25" = 20.84514502451109
    "O25" = 19.39872123832938
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
